# "Version 1." -> "Version 2."
# The target markup splits "Version" into two runs ("Versi" / "on"),
# drops the trailing "." from the " 1." run (leaving " 2"), and adds a
# brand-new run containing just "." placed after the (unchanged)
# _GoBack bookmark. Plain Find/Replace normalizes/merges runs, so the
# run split is produced with Range.InsertXML (which inserts literal
# OOXML without Word's usual "merge same-formatted runs" clean-up),
# while the rest is done with ordinary Range.Text / Delete / InsertAfter.

$d = $word.ActiveDocument

# --- Step 1: split "Version" (chars 0-6) into "Versi" + "on" -------------
# Replacing only the "on" substring (chars 5-6) with an explicit <w:r>
# leaves the existing "Versi" run untouched and keeps the spellEnd
# proofErr marker right after the new "on" run, matching the target.
$rOn = $d.Range(5, 7)
$onXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'
$rOn.InsertXML($onXml)

# --- Step 2: "1" -> "2" (still inside the " 1." run; same formatting so ---
# this is a plain in-place text edit, no run split needed)
$rDigit = $d.Range(8, 9)
$rDigit.Text = "2"

# --- Step 3: drop the trailing "." from that run (it reappears after the -
# bookmark as its own run in step 4)
$rPeriod = $d.Range(9, 10)
$rPeriod.Delete()

# --- Step 4: insert a fresh "." run right after the _GoBack bookmark -----
$bm = $d.Bookmarks("_GoBack")
$rEnd = $d.Range($bm.End, $bm.End)
$rEnd.InsertAfter(".")
